# Apply the "Add files via upload" update to CORE_holdings.xlsx:
#  1. Bump the "as of" date in the confidential disclaimer text (A11) from
#     2021-03-23 to 2021-03-24.
#  2. Refresh the Weight (col D) and Percent Change (col E) figures for
#     rows 2-8 with the new model-holdings numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; temporarily unprotect so the cells can be
# edited, then restore protection afterwards.
$ws.Unprotect()

# --- Disclaimer text: update the "as of" date -----------------------------
$oldText = $ws.Range("A11").Value2
$newText = $oldText -replace "2021-03-23", "2021-03-24"
$ws.Range("A11").Value = $newText

# --- Updated Weight / Percent Change figures -------------------------------
$ws.Range("D2").Value = 0.4971088642703818
$ws.Range("E2").Value = 0.0009367343997694899

$ws.Range("D3").Value = 0.2465232797940843
$ws.Range("E3").Value = -0.01117491851621932

$ws.Range("D4").Value = 0.09780548032580662
$ws.Range("E4").Value = -0.01035658101730474

$ws.Range("D5").Value = 0.1004993012062562
$ws.Range("E5").Value = -0.001620089104900746

$ws.Range("D6").Value = 0.0299421976376355
$ws.Range("E6").Value = -0.01223491027732482

$ws.Range("D7").Value = 0.02812087676583561
$ws.Range("E7").Value = -0.02189316137250219

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = -0.004446961788895698

# Restore the original sheet protection.
$ws.Protect("D382")

Write-Output "Updated disclaimer date and D2:E8 weight/percent-change figures."
